# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# values for the file "d0017418-62ea-4103-aa37-62df6560000a.md" after a new
# handoff xliff generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-12 12:48:15"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-12 12:48:08"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-12 12:48:15"
